$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.039365724995899
$ws.Range("D2").Value = 1.049414599242759
$ws.Range("E2").Value = 1.04777698544212
$ws.Range("F2").Value = 1.058352371740531
$ws.Range("I2").Value = 1.035471182256469
$ws.Range("J2").Value = 1.044458133295028
$ws.Range("K2").Value = 1.052171723055868
$ws.Range("L2").Value = 1.050538676133502
$ws.Range("M2").Value = 1.061084840223917
$ws.Range("N2").Value = 1.01874652582255
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.040305033722873
$ws.Range("D3").Value = 1.050297490872406
$ws.Range("E3").Value = 1.048623427870858
$ws.Range("F3").Value = 1.059310328472965
$ws.Range("I3").Value = 1.035598168789212
$ws.Range("J3").Value = 1.04504266888431
$ws.Range("K3").Value = 1.052866673585557
$ws.Range("L3").Value = 1.051196941466026
$ws.Range("M3").Value = 1.061856448149011
$ws.Range("N3").Value = 1.018940185538238
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04091285894864
$ws.Range("D4").Value = 1.050869139592445
$ws.Range("E4").Value = 1.049171584786451
$ws.Range("F4").Value = 1.05993081325743
$ws.Range("I4").Value = 1.035678156644725
$ws.Range("J4").Value = 1.045420307439825
$ws.Range("K4").Value = 1.053316074325445
$ws.Range("L4").Value = 1.051622689984296
$ws.Range("M4").Value = 1.062355728618046
$ws.Range("N4").Value = 1.019065279602891
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.041168394757145
$ws.Range("D5").Value = 1.051109545239993
$ws.Range("E5").Value = 1.049402136812235
$ws.Range("F5").Value = 1.060191812762787
$ws.Range("I5").Value = 1.035711260589304
$ws.Range("J5").Value = 1.045578923177098
$ws.Range("K5").Value = 1.053504934498214
$ws.Range("L5").Value = 1.051801627297735
$ws.Range("M5").Value = 1.062565624549014
$ws.Range("N5").Value = 1.019117816887041
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.041211300670829
$ws.Range("D6").Value = 1.051149915322668
$ws.Range("E6").Value = 1.049440853738211
$ws.Range("F6").Value = 1.060235644349295
$ws.Range("I6").Value = 1.035716788196843
$ws.Range("J6").Value = 1.045605547010981
$ws.Range("K6").Value = 1.053536640936464
$ws.Range("L6").Value = 1.051831668849948
$ws.Range("M6").Value = 1.062600866870974
$ws.Range("N6").Value = 1.019126635052622
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040916273408002
$ws.Range("D7").Value = 1.05087235157408
$ws.Range("E7").Value = 1.049174665015751
$ws.Range("F7").Value = 1.059934300166969
$ws.Range("I7").Value = 1.03567860103715
$ws.Range("J7").Value = 1.045422427436308
$ws.Range("K7").Value = 1.053318598152595
$ws.Range("L7").Value = 1.05162508114096
$ws.Range("M7").Value = 1.062358533265507
$ws.Range("N7").Value = 1.01906598181448
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.039683162023555
$ws.Range("D8").Value = 1.049712901488923
$ws.Range("E8").Value = 1.048062950201645
$ws.Range("F8").Value = 1.058675987937798
$ws.Range("I8").Value = 1.035514549026451
$ws.Range("J8").Value = 1.044655802413753
$ws.Range("K8").Value = 1.052406641850534
$ws.Range("L8").Value = 1.050761179311078
$ws.Range("M8").Value = 1.061345608298706
$ws.Range("N8").Value = 1.01881201867774
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.037510533938304
$ws.Range("D9").Value = 1.047672610336433
$ws.Range("E9").Value = 1.046107488833
$ws.Range("F9").Value = 1.056463510749494
$ws.Range("I9").Value = 1.035208801414458
$ws.Range("J9").Value = 1.043300399780735
$ws.Range("K9").Value = 1.050797574888201
$ws.Range("L9").Value = 1.049237443421727
$ws.Range("M9").Value = 1.059560747793674
$ws.Range("N9").Value = 1.01836285985382
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.03606235574764
$ws.Range("D10").Value = 1.046314383593741
$ws.Range("E10").Value = 1.04480629537158
$ws.Range("F10").Value = 1.054991853009912
$ws.Range("I10").Value = 1.03499381156619
$ws.Range("J10").Value = 1.042393825469102
$ws.Range("K10").Value = 1.049723530387952
$ws.Range("L10").Value = 1.04822072236117
$ws.Range("M10").Value = 1.058370940573601
$ws.Range("N10").Value = 1.018062337543619
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035435347172047
$ws.Range("D11").Value = 1.045726739301391
$ws.Range("E11").Value = 1.044243461035725
$ws.Range("F11").Value = 1.05435541441579
$ws.Range("I11").Value = 1.034898081793185
$ws.Range("J11").Value = 1.042000576078215
$ws.Range("K11").Value = 1.049258154835188
$ws.Range("L11").Value = 1.047780271338031
$ws.Range("M11").Value = 1.05785577909278
$ws.Range("N11").Value = 1.017931955622681
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.035202458454845
$ws.Range("D12").Value = 1.045508534947021
$ws.Range("E12").Value = 1.044034489367494
$ws.Range("F12").Value = 1.054119134077717
$ws.Range("I12").Value = 1.03486212793302
$ws.Range("J12").Value = 1.041854402009887
$ws.Range("K12").Value = 1.049085248538905
$ws.Range("L12").Value = 1.047616638570509
$ws.Range("M12").Value = 1.057664431241155
$ws.Range("N12").Value = 1.017883488186244
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035252413400738
$ws.Range("D13").Value = 1.045555337202693
$ws.Range("E13").Value = 1.044079310403778
$ws.Range("F13").Value = 1.054169811503294
$ws.Range("I13").Value = 1.034869858046338
$ws.Range("J13").Value = 1.04188576152294
$ws.Range("K13").Value = 1.049122339550143
$ws.Range("L13").Value = 1.047651739679919
$ws.Range("M13").Value = 1.057705475712777
$ws.Range("N13").Value = 1.017893886320314
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035416096305354
$ws.Range("D14").Value = 1.04570870096698
$ws.Range("E14").Value = 1.044226185527834
$ws.Range("F14").Value = 1.0543358809277
$ws.Range("I14").Value = 1.034895117900991
$ws.Range("J14").Value = 1.041988495394875
$ws.Range("K14").Value = 1.049243863260634
$ws.Range("L14").Value = 1.047766746001345
$ws.Range("M14").Value = 1.057839962089613
$ws.Range("N14").Value = 1.017927950058118
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035516948120115
$ws.Range("D15").Value = 1.045803203141397
$ws.Range("E15").Value = 1.04431669211153
$ws.Range("F15").Value = 1.054438217891313
$ws.Range("I15").Value = 1.034910628937322
$ws.Range("J15").Value = 1.042051779406139
$ws.Range("K15").Value = 1.049318732085434
$ws.Range("L15").Value = 1.047837601293519
$ws.Range("M15").Value = 1.057922824421308
$ws.Range("N15").Value = 1.017948932854418
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036103969623606
$ws.Range("D16").Value = 1.046353393764406
$ws.Range("E16").Value = 1.044843661393491
$ws.Range("F16").Value = 1.05503410828374
$ws.Range("I16").Value = 1.035000109324263
$ws.Range("J16").Value = 1.042419909537178
$ws.Range("K16").Value = 1.049754409448871
$ws.Range("L16").Value = 1.048249949431855
$ws.Range("M16").Value = 1.05840513093308
$ws.Range("N16").Value = 1.01807098524471
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036472209941274
$ws.Range("D17").Value = 1.046698642364059
$ws.Range("E17").Value = 1.045174374433453
$ws.Range("F17").Value = 1.055408109371809
$ws.Range("I17").Value = 1.035055532325885
$ws.Range("J17").Value = 1.04265064197418
$ws.Range("K17").Value = 1.050027616852184
$ws.Range("L17").Value = 1.048508550602672
$ws.Range("M17").Value = 1.058707678664073
$ws.Range("N17").Value = 1.018147477782951
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.036687004227081
$ws.Range("D18").Value = 1.046900065956187
$ws.Range("E18").Value = 1.04536733069447
$ws.Range("F18").Value = 1.055626334878503
$ws.Range("I18").Value = 1.035087605071824
$ws.Range("J18").Value = 1.042785157013816
$ws.Range("K18").Value = 1.050186944397694
$ws.Range("L18").Value = 1.048659368530462
$ws.Range("M18").Value = 1.058884152736385
$ws.Range("N18").Value = 1.01819207007775
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03676024450913
$ws.Range("D19").Value = 1.04696875394712
$ws.Range("E19").Value = 1.045433133412303
$ws.Range("F19").Value = 1.055700757158727
$ws.Range("I19").Value = 1.035098497850312
$ws.Range("J19").Value = 1.042831011724347
$ws.Range("K19").Value = 1.050241265901614
$ws.Range("L19").Value = 1.0487107901637
$ws.Range("M19").Value = 1.058944326326178
$ws.Range("N19").Value = 1.018207270730398
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036432700619693
$ws.Range("D20").Value = 1.046661595697307
$ws.Range("E20").Value = 1.045138886155294
$ws.Range("F20").Value = 1.055367974620936
$ws.Range("I20").Value = 1.035049612286628
$ws.Range("J20").Value = 1.042625893520292
$ws.Range("K20").Value = 1.049998307348815
$ws.Range("L20").Value = 1.048480807180887
$ws.Range("M20").Value = 1.058675217827086
$ws.Range("N20").Value = 1.01813927339007
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035367895481666
$ws.Range("D21").Value = 1.045663537138396
$ws.Range("E21").Value = 1.044182931968183
$ws.Range("F21").Value = 1.054286974262973
$ws.Range("I21").Value = 1.03488769041694
$ws.Range("J21").Value = 1.04195824568624
$ws.Range("K21").Value = 1.04920807879318
$ws.Range("L21").Value = 1.047732880316399
$ws.Range("M21").Value = 1.057800359027395
$ws.Range("N21").Value = 1.017917920176692
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034698469576369
$ws.Range("D22").Value = 1.04503644050876
$ws.Range("E22").Value = 1.043582407337261
$ws.Range("F22").Value = 1.053608008819879
$ws.Range("I22").Value = 1.034783595478347
$ws.Range("J22").Value = 1.041537869275632
$ws.Range("K22").Value = 1.048710970348338
$ws.Range("L22").Value = 1.047262457985065
$ws.Range("M22").Value = 1.05725033618638
$ws.Range("N22").Value = 1.017778528176612
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035053339006644
$ws.Range("D23").Value = 1.045368835839978
$ws.Range("E23").Value = 1.043900706932134
$ws.Range("F23").Value = 1.053967874248908
$ws.Range("I23").Value = 1.034838994830628
$ws.Range("J23").Value = 1.041760775298662
$ws.Range("K23").Value = 1.048974521270928
$ws.Range("L23").Value = 1.047511853637171
$ws.Range("M23").Value = 1.057541909954375
$ws.Range("N23").Value = 1.017852443106324
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036450553163117
$ws.Range("D24").Value = 1.046678335349546
$ws.Range("E24").Value = 1.045154921604784
$ws.Range("F24").Value = 1.055386109549447
$ws.Range("I24").Value = 1.035052288084118
$ws.Range("J24").Value = 1.04263707648928
$ws.Range("K24").Value = 1.050011551143643
$ws.Range("L24").Value = 1.048493343300112
$ws.Range("M24").Value = 1.058689885472556
$ws.Range("N24").Value = 1.018142980677809
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03807217219489
$ws.Range("D25").Value = 1.048199733527884
$ws.Range("E25").Value = 1.046612596869753
$ws.Range("F25").Value = 1.057034908536952
$ws.Range("I25").Value = 1.035289814662173
$ws.Range("J25").Value = 1.04365133219559
$ws.Range("K25").Value = 1.051213797171176
$ws.Range("L25").Value = 1.049631528505565
$ws.Range("M25").Value = 1.060022164605759
$ws.Range("N25").Value = 1.01847917066416
